# res_bus/vm_pu.xlsx - "case with 380 kV done"
# The slack-bus voltage setpoint (column B) was changed from 1.05 p.u. to
# 1.02 p.u., and the power-flow was re-run, which changed the resulting
# per-unit voltage magnitudes for every bus (columns C-F and I-N) in every
# row of the results table (rows 2-25). Column A (bus index), column G
# (always 1) and column H (always blank) are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027081630696739
$ws.Range("D2").Value = 1.031976431707556
$ws.Range("E2").Value = 1.051185743070378
$ws.Range("F2").Value = 1.056252872354429
$ws.Range("I2").Value = 1.036484085736854
$ws.Range("J2").Value = 1.032241487854122
$ws.Range("K2").Value = 1.034782994692483
$ws.Range("L2").Value = 1.053937944986846
$ws.Range("M2").Value = 1.058991091804023
$ws.Range("N2").Value = 1.014789840171217

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027819378218896
$ws.Range("D3").Value = 1.032516163129518
$ws.Range("E3").Value = 1.052201094046537
$ws.Range("F3").Value = 1.05727489483773
$ws.Range("I3").Value = 1.036658014836282
$ws.Range("J3").Value = 1.032620357739459
$ws.Range("K3").Value = 1.035132122626628
$ws.Range("L3").Value = 1.054765370111054
$ws.Range("M3").Value = 1.059826185945831
$ws.Range("N3").Value = 1.014915043249939

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.028297121612007
$ws.Range("D4").Value = 1.032865646051346
$ws.Range("E4").Value = 1.052859331786878
$ws.Range("F4").Value = 1.057937252715385
$ws.Range("I4").Value = 1.036769379671436
$ws.Range("J4").Value = 1.0328651907583
$ws.Range("K4").Value = 1.03535756313696
$ws.Range("L4").Value = 1.055301395536222
$ws.Range("M4").Value = 1.06036697753495
$ws.Range("N4").Value = 1.014995944155833

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028498052129278
$ws.Range("D5").Value = 1.033012624550583
$ws.Range("E5").Value = 1.053136350057217
$ws.Range("F5").Value = 1.058215955884468
$ws.Range("I5").Value = 1.036815914449214
$ws.Range("J5").Value = 1.032968040560766
$ws.Range("K5").Value = 1.035452224986213
$ws.Range("L5").Value = 1.055526889431897
$ws.Range("M5").Value = 1.060594428036845
$ws.Range("N5").Value = 1.015029927248775

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028531794310324
$ws.Range("D6").Value = 1.033037306113961
$ws.Range("E6").Value = 1.053182879915242
$ws.Range("F6").Value = 1.058262765873286
$ws.Range("I6").Value = 1.036823711225815
$ws.Range("J6").Value = 1.032985304888719
$ws.Range("K6").Value = 1.035468112447218
$ws.Range("L6").Value = 1.055564759544347
$ws.Range("M6").Value = 1.060632623907335
$ws.Range("N6").Value = 1.015035631526975

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028299806114956
$ws.Range("D7").Value = 1.032867609767851
$ws.Range("E7").Value = 1.05286303216042
$ws.Range("F7").Value = 1.057940975788243
$ws.Range("I7").Value = 1.036770002584359
$ws.Range("J7").Value = 1.032866565350281
$ws.Range("K7").Value = 1.035358828460029
$ws.Range("L7").Value = 1.055304408013314
$ws.Range("M7").Value = 1.060370016341488
$ws.Range("N7").Value = 1.014996398348712

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027330877986616
$ws.Range("D8").Value = 1.032158785616958
$ws.Range("E8").Value = 1.051528628676912
$ws.Range("F8").Value = 1.056598053311973
$ws.Range("I8").Value = 1.036543109535508
$ws.Range("J8").Value = 1.032369594525596
$ws.Range("K8").Value = 1.034901080446795
$ws.Range("L8").Value = 1.054217446911642
$ws.Range("M8").Value = 1.059273226063735
$ws.Range("N8").Value = 1.014832176478643

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.025626423547644
$ws.Range("D9").Value = 1.030911665571131
$ws.Range("E9").Value = 1.049186772138971
$ws.Range("F9").Value = 1.054239682531164
$ws.Range("I9").Value = 1.036134299557895
$ws.Range("J9").Value = 1.031491458371867
$ws.Range("K9").Value = 1.034090931986845
$ws.Range("L9").Value = 1.052306924625686
$ws.Range("M9").Value = 1.057343886025878
$ws.Range("N9").Value = 1.014541941759744

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024492186793902
$ws.Range("D10").Value = 1.030081644689576
$ws.Range("E10").Value = 1.047632014581135
$ws.Range("F10").Value = 1.05267291103865
$ws.Range("I10").Value = 1.035855756887132
$ws.Range("J10").Value = 1.03090448040382
$ws.Range("K10").Value = 1.033548522980391
$ws.Range("L10").Value = 1.051036561722458
$ws.Range("M10").Value = 1.056059978149116
$ws.Range("N10").Value = 1.014347900407533

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024001561398558
$ws.Range("D11").Value = 1.029722588409846
$ws.Range("E11").Value = 1.046960339044985
$ws.Range("F11").Value = 1.051995796335812
$ws.Range("I11").Value = 1.035733730592286
$ws.Range("J11").Value = 1.030649957903153
$ws.Range("K11").Value = 1.033313120448971
$ws.Range("L11").Value = 1.050487280764635
$ws.Range("M11").Value = 1.055504596412129
$ws.Range("N11").Value = 1.014263752421027

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023819399125641
$ws.Range("D12").Value = 1.0295892729924
$ws.Range("E12").Value = 1.046711082076628
$ws.Range("F12").Value = 1.05174448350173
$ws.Range("I12").Value = 1.035688192679507
$ws.Range("J12").Value = 1.030555364372796
$ws.Range("K12").Value = 1.033225602119187
$ws.Range("L12").Value = 1.050283373757654
$ws.Range("M12").Value = 1.055298388008442
$ws.Range("N12").Value = 1.014232477452472

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023858470001129
$ws.Range("D13").Value = 1.02961786712704
$ws.Range("E13").Value = 1.046764537969669
$ws.Range("F13").Value = 1.051798381990604
$ws.Range("I13").Value = 1.035697970293932
$ws.Range("J13").Value = 1.030575657365797
$ws.Range("K13").Value = 1.03324437867013
$ws.Range("L13").Value = 1.050327107038104
$ws.Range("M13").Value = 1.055342616541608
$ws.Range("N13").Value = 1.014239186877517

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023986502210239
$ws.Range("D14").Value = 1.029711567402686
$ws.Range("E14").Value = 1.046939730621678
$ws.Range("F14").Value = 1.051975018697976
$ws.Range("I14").Value = 1.035729970732866
$ws.Range("J14").Value = 1.030642139837777
$ws.Range("K14").Value = 1.033305887768524
$ws.Range("L14").Value = 1.050470423285818
$ws.Range("M14").Value = 1.055487549418482
$ws.Range("N14").Value = 1.014261167600765

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024065397446773
$ws.Range("D15").Value = 1.02976930645239
$ws.Range("E15").Value = 1.047047703536736
$ws.Range("F15").Value = 1.052083876639767
$ws.Range("I15").Value = 1.035749659201159
$ws.Range("J15").Value = 1.030683094954229
$ws.Range("K15").Value = 1.033343775074268
$ws.Range("L15").Value = 1.050558741140611
$ws.Range("M15").Value = 1.05557685865836
$ws.Range("N15").Value = 1.014274708188899

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024524758807009
$ws.Range("D16").Value = 1.030105481559078
$ws.Range("E16").Value = 1.047676624118539
$ws.Range("F16").Value = 1.052717876569519
$ws.Range("I16").Value = 1.03586382561158
$ws.Range("J16").Value = 1.030921364794414
$ws.Range("K16").Value = 1.03356413468017
$ws.Range("L16").Value = 1.051073032518374
$ws.Range("M16").Value = 1.056096848877461
$ws.Range("N16").Value = 1.014353482394084

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024813040907517
$ws.Range("D17").Value = 1.030316449794441
$ws.Range("E17").Value = 1.048071543767124
$ws.Range("F17").Value = 1.05311591896182
$ws.Range("I17").Value = 1.035935060746958
$ws.Range("J17").Value = 1.031070730336985
$ws.Range("K17").Value = 1.033702217793935
$ws.Range("L17").Value = 1.051395847360066
$ws.Range("M17").Value = 1.056423175427204
$ws.Range("N17").Value = 1.014402861699712

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024981239743167
$ws.Range("D18").Value = 1.030439537479612
$ws.Range("E18").Value = 1.048302042684701
$ws.Range("F18").Value = 1.053348216444191
$ws.Range("I18").Value = 1.035976474372822
$ws.Range("J18").Value = 1.031157818194475
$ws.Range("K18").Value = 1.033782707558527
$ws.Range("L18").Value = 1.051584216390453
$ws.Range("M18").Value = 1.056613570026381
$ws.Range("N18").Value = 1.014431651552682

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025038599417941
$ws.Range("D19").Value = 1.030481512825241
$ws.Range("E19").Value = 1.048380662091594
$ws.Range("F19").Value = 1.05342744523159
$ws.Range("I19").Value = 1.035990572152623
$ws.Range("J19").Value = 1.031187507011276
$ws.Range("K19").Value = 1.033810143670269
$ws.Range("L19").Value = 1.051648458322675
$ws.Range("M19").Value = 1.056678498762639
$ws.Range("N19").Value = 1.014441466057046

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024782105910035
$ws.Range("D20").Value = 1.030293811423435
$ws.Range("E20").Value = 1.048029157209512
$ws.Range("F20").Value = 1.053073199729115
$ws.Range("I20").Value = 1.035927432017307
$ws.Range("J20").Value = 1.031054708392031
$ws.Range("K20").Value = 1.033687408130072
$ws.Range("L20").Value = 1.051361204452671
$ws.Range("M20").Value = 1.056388158099172
$ws.Range("N20").Value = 1.014397565032887

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023948797752989
$ws.Range("D21").Value = 1.029683973500229
$ws.Range("E21").Value = 1.046888134304598
$ws.Range("F21").Value = 1.051922998122244
$ws.Range("I21").Value = 1.035720553245048
$ws.Range("J21").Value = 1.030622563848218
$ws.Range("K21").Value = 1.033287777053237
$ws.Range("L21").Value = 1.050428216891124
$ws.Range("M21").Value = 1.05544486794018
$ws.Range("N21").Value = 1.014254695337588

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02342531507955
$ws.Range("D22").Value = 1.029300857900241
$ws.Range("E22").Value = 1.046172077900707
$ws.Range("F22").Value = 1.051200965261707
$ws.Range("I22").Value = 1.035589254652746
$ws.Range("J22").Value = 1.030350553988877
$ws.Range("K22").Value = 1.033036054620593
$ws.Range("L22").Value = 1.049842307674263
$ws.Range("M22").Value = 1.054852276919702
$ws.Range("N22").Value = 1.01416475972342

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023702779798712
$ws.Range("D23").Value = 1.029503924408904
$ws.Range("E23").Value = 1.04655154459976
$ws.Range("F23").Value = 1.051583619634525
$ws.Range("I23").Value = 1.035658974453361
$ws.Range("J23").Value = 1.030494779950431
$ws.Range("K23").Value = 1.033169540520666
$ws.Range("L23").Value = 1.050152842816373
$ws.Range("M23").Value = 1.055166373524204
$ws.Range("N23").Value = 1.014212446378854

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024796083953358
$ws.Range("D24").Value = 1.030304040625606
$ws.Range("E24").Value = 1.048048309410199
$ws.Range("F24").Value = 1.053092502322059
$ws.Range("I24").Value = 1.035930879534091
$ws.Range("J24").Value = 1.031061948126147
$ws.Range("K24").Value = 1.03369410013983
$ws.Range("L24").Value = 1.051376857854756
$ws.Range("M24").Value = 1.056403980756137
$ws.Range("N24").Value = 1.014399958406564

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.026066709194281
$ws.Range("D25").Value = 1.031233838046156
$ws.Range("E25").Value = 1.049791061221498
$ws.Range("F25").Value = 1.05484841841079
$ws.Range("I25").Value = 1.03624104807203
$ws.Range("J25").Value = 1.031718756361695
$ws.Range("K25").Value = 1.034300786965223
$ws.Range("L25").Value = 1.052800260140162
$ws.Range("M25").Value = 1.057842263116105
$ws.Range("N25").Value = 1.014617073298836
